# Updates cryptos list values (price + 1h volume%), generated from the
# GitHub Actions data refresh on Fri Mar 17 20:27:35 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.759.66"
$ws.Range("E2").Value = "  +7.32%  "

# Row 3
$ws.Range("D3").Value = "1.744.90"
$ws.Range("E3").Value = "  +4.01%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.27"
$ws.Range("E5").Value = "  +2.07%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9977"
$ws.Range("E6").Value = "  -0.21%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.38"
$ws.Range("E8").Value = "  +2.91%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3392"
$ws.Range("E9").Value = "  +4.27%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.192"
$ws.Range("E10").Value = "  +3.98%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07487"
$ws.Range("E11").Value = "  +5.68%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9967"
$ws.Range("E12").Value = "  -0.40%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.422"
$ws.Range("E13").Value = "  +5.33%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.57"
$ws.Range("E14").Value = "  +4.67%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.080"
$ws.Range("E15").Value = "  +6.57%  "

# Row 16
$ws.Range("D16").Value = "1.745.03"
$ws.Range("E16").Value = "  +3.85%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001080"
$ws.Range("E17").Value = "  +2.94%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06739"
$ws.Range("E18").Value = "  +2.26%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.89"
$ws.Range("E19").Value = "  +5.03%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9985"
$ws.Range("E20").Value = "  -0.14%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.78"
$ws.Range("E21").Value = "  +5.40%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.227"
$ws.Range("E22").Value = "  +5.06%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.82"
$ws.Range("E23").Value = "  -0.47%  "

# Row 24
$ws.Range("D24").Value = "26.776.74"
$ws.Range("E24").Value = "  +7.33%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.469"
$ws.Range("E25").Value = "  +0.76%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.474"
$ws.Range("E26").Value = "  +24.40%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.442"
$ws.Range("E27").Value = "  +1.14%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.96"
$ws.Range("E28").Value = "  +2.60%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.67"
$ws.Range("E29").Value = "  +4.98%  "

# Row 30
$ws.Range("D30").Value = "1.938.94"
$ws.Range("E30").Value = "  +4.03%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "132.79"
$ws.Range("E31").Value = "  +5.51%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.120"
$ws.Range("E32").Value = "  +1.13%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.072"
$ws.Range("E33").Value = "  +5.14%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08609"
$ws.Range("E34").Value = "  +1.28%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.693"
$ws.Range("E35").Value = "  +3.02%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.94"
$ws.Range("E36").Value = "  +5.19%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.444"
$ws.Range("E37").Value = "  +5.01%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02354"
$ws.Range("E38").Value = "  +4.48%  "

# Row 39
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2183"
$ws.Range("E39").Value = "  +4.20%  "

# Row 40
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06274"
$ws.Range("E40").Value = "  +4.33%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.506"
$ws.Range("E41").Value = "  +3.30%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.227"
$ws.Range("E42").Value = "  -0.38%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6276"
$ws.Range("E43").Value = "  +5.34%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.38"
$ws.Range("E44").Value = "  +5.07%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9973"
$ws.Range("E45").Value = "  -0.20%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.937"
$ws.Range("E46").Value = "  +2.45%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6098"
$ws.Range("E47").Value = "  +6.42%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.45"
$ws.Range("E48").Value = "  +3.13%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.075"
$ws.Range("E49").Value = "  +5.53%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07231"
$ws.Range("E50").Value = "  +2.91%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.84"
$ws.Range("E51").Value = "  +4.04%  "
